$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.543.84"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "1.848.77"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.032"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +2.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.41"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.026"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  +2.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3777"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07384"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8736"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.54%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").Value = "1.860.89"
$ws.Range("E12").Value = "  -8.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.511"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.678"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07193"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.92"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.032"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009033"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("E19").Value = "  +2.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "27.563.03"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.248"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.33"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.60"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.51%  "
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.69"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.969"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.269"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.90"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09043"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.194"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7592"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.509"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.875"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01969"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05291"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.816"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5141"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1672"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.738"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.457"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.51"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.56"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.705"
$ws.Range("D46").ClearFormats()
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06395"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4634"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.844"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.15"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.87"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.28%  "
